$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Row 25: "Agregar donde ver el nuevo proveedor de impreso" marked done (si),
#     with start/finish dates, then hidden (matches other completed rows) ---
$ws.Cells.Item(5,4).Copy()
$ws.Cells.Item(25,4).PasteSpecial(-4122)
$ws.Cells.Item(25,5).PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Cells.Item(25,3).Value = "si"
$ws.Cells.Item(25,4).Value = 43805
$ws.Cells.Item(25,5).Value = 43805

# --- Row 32: "preguntar si afectara a la formula de pyg" sub-task gets flagged (red font) ---
$ws.Cells.Item(32,2).Font.Color = 255

# --- Row 35: "Agregar estatus a pedido" marked done (si), with dates, then hidden ---
$ws.Cells.Item(5,4).Copy()
$ws.Cells.Item(35,4).PasteSpecial(-4122)
$ws.Cells.Item(35,5).PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Cells.Item(35,3).Value = "si"
$ws.Cells.Item(35,4).Value = 43804
$ws.Cells.Item(35,5).Value = 43804

# --- Row 42 & 44: clear the red-flag font (task resolved) ---
$ws.Cells.Item(42,2).Font.Color = 0
$ws.Cells.Item(44,2).Font.Color = 0

# --- Row 46: add the missing "Fecha fin" to match "Fecha inicio" ---
$ws.Rows.Item(46).Hidden = $false
$ws.Cells.Item(5,4).Copy()
$ws.Cells.Item(46,5).PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Cells.Item(46,5).Value = 43804
$ws.Rows.Item(46).Hidden = $true

# --- Row 49: flag with red font (pending task) ---
$ws.Cells.Item(49,2).Font.Color = 255

# Now that all direct-formatting copies into rows 25/35 are done, restore their
# hidden state to match the other completed (si) rows in the table.
$ws.Rows.Item(25).Hidden = $true
$ws.Rows.Item(35).Hidden = $true

# --- New tasks appended to the list ---
$ws.Cells.Item(52,1).Value = 19
$ws.Cells.Item(52,2).Value = "Acomodar todos los rs y st"
$ws.Cells.Item(52,3).Value = "no"

$ws.Cells.Item(53,1).Value = 20
$ws.Cells.Item(53,2).Value = "si el pedido ya se pago, que se muestre 0 en resto"
$ws.Cells.Item(53,3).Value = "no"

$ws.Cells.Item(54,1).Value = 21
$ws.Cells.Item(54,2).Value = "Preguntar que hacer con sticky"
$ws.Cells.Item(54,2).Font.Color = 255
$ws.Cells.Item(54,3).Value = "no"

$ws.Cells.Item(55,1).Value = 22
$ws.Cells.Item(55,2).Value = "copia en dos pc"
$ws.Cells.Item(55,2).Font.Color = 255
$ws.Cells.Item(55,3).Value = "no"

$ws.Cells.Item(56,1).Value = 23
$ws.Cells.Item(56,2).Value = "copia en la nube"
$ws.Cells.Item(56,3).Value = "no"

# Update the worksheet selection to match where editing ended.
$ws.Range("B55").Select()

# Grow the table / autofilter range to cover the newly added rows.
$tbl = $ws.ListObjects.Item("Tabla2")
$tbl.Resize($ws.Range("A1:E56"))
